# Terra Machina Reference Guide - update
# - Bestiary: rename "ElementalWeaknesses" column to "Weakness" and fill in
#   two new monsters (Muskroom, Flying Drone) into the pre-existing blank
#   table rows.
# - Items: add a new "Dangerous Mushroom" item row.
# - Selections / active sheet tidy-up left by the author's last save.

$wb = $excel.ActiveWorkbook

$items    = $wb.Worksheets.Item("Items")
$bestiary = $wb.Worksheets.Item("Bestiary")
$moves    = $wb.Worksheets.Item("SampleCustomMoves")

# ---------------------------------------------------------------------
# Bestiary sheet: rename the "ElementalWeaknesses" header to "Weakness"
# (this also renames the Table2 column automatically) and populate the
# two new creature rows that already existed as blank table rows.
# ---------------------------------------------------------------------

$bestiary.Range("A4").Value = "Muskroom"
$bestiary.Range("B4").Value = "A muskrat covered in mushrooms"
$bestiary.Range("C4").Value = 3
$bestiary.Range("D4").Value = 0
$bestiary.Range("E4").Value = "Fire"
$bestiary.Range("F4").Value = "Scratch: 1d4+3 damage`r`nBite: 1d4 damage.  Causes Poison.`r`nSpore Cloud: Cause Poison on all enemies"
$bestiary.Rows.Item(4).RowHeight = 75

$bestiary.Range("E1").Value = "Weakness"

$bestiary.Range("A5").Value = "Flying Drone"
$bestiary.Range("B5").Value = "A robot in the air that can shoot an enemy"
$bestiary.Range("C5").Value = 10
$bestiary.Range("D5").Value = 1
$bestiary.Range("E5").Value = "Lightning"
$bestiary.Range("F5").Value = "Laser Shot: 1d6 damage"

# ---------------------------------------------------------------------
# Items sheet: add the new "Dangerous Mushroom" item.
# ---------------------------------------------------------------------

$items.Range("A14").Value = "Dangerous Mushroom"
$items.Range("B14").Value = "Requires a resilience roll:`r`n-`tSuccess: Heal 1d6 HP`r`n-`tMinor Success: Heal 1d6 HP but gain Poison.`r`n-`tFail: Gain Poison"
$items.Range("B14").WrapText = $true
$items.Range("C14").Value = 5
$items.Rows.Item(14).RowHeight = 60

$items.Columns.Item(1).ColumnWidth = 22.8

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping to match the author's last save.
# ---------------------------------------------------------------------

$bestiary.Range("F3").Select() | Out-Null
$moves.Range("B3").Select() | Out-Null
$items.Range("C15").Select() | Out-Null
